$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the query text in B2: remove the tp.endocrine_therapy_type / head(labels(samp)) clauses
$oldQuery = $ws.Range("B2").Value2
$newQuery = $oldQuery.Replace('and tp.endocrine_therapy_type IN ["Other"]  and head(labels(samp)) IN ["sample"]', ' ')
$ws.Range("B2").Value = $newQuery

# Adjust row height for row 2
$ws.Rows.Item(2).RowHeight = 375

# Adjust the view: scroll so row 2 is at the top, and move the selection to B2
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
